$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(100, 8).Value = 2548.24
$ws.Cells.Item(100, 9).Value = 1820.8
$ws.Cells.Item(100, 10).Value = 3033.2
$ws.Cells.Item(100, 11).Value = 1820.8
$ws.Cells.Item(100, 12).Value = 3033.2
$ws.Cells.Item(100, 13).Value = -1279.8
$ws.Cells.Item(100, 14).Value = -4115.2

$ws.Cells.Item(107, 8).Value = 29004.229
$ws.Cells.Item(107, 9).Value = 40384.36
$ws.Cells.Item(107, 10).Value = 553.9
$ws.Cells.Item(107, 11).Value = 40384.36
$ws.Cells.Item(107, 12).Value = 553.9
$ws.Cells.Item(107, 13).Value = -38464.36
$ws.Cells.Item(107, 14).Value = -4393.9

$ws.Cells.Item(112, 8).Value = 1022.89
$ws.Cells.Item(112, 9).Value = 500
$ws.Cells.Item(112, 10).Value = 1028.1718
$ws.Cells.Item(112, 11).Value = 1500
$ws.Cells.Item(112, 12).Value = 3084.5154
$ws.Cells.Item(112, 13).Value = -392
$ws.Cells.Item(112, 14).Value = -5300.5154

$ws.Cells.Item(137, 8).Value = 1669.0303
$ws.Cells.Item(137, 9).Value = 1328.8422
$ws.Cells.Item(137, 10).Value = 2130.7144
$ws.Cells.Item(137, 11).Value = 3986.5266
$ws.Cells.Item(137, 12).Value = 6392.1432
$ws.Cells.Item(137, 13).Value = -1436.5266
$ws.Cells.Item(137, 14).Value = -11492.1432

$ws.Cells.Item(138, 8).Value = 3917.7273
$ws.Cells.Item(138, 9).Value = 1174.0294
$ws.Cells.Item(138, 10).Value = 8359.904
$ws.Cells.Item(138, 11).Value = 3522.0882
$ws.Cells.Item(138, 12).Value = 25079.712
$ws.Cells.Item(138, 13).Value = 1617.9118
$ws.Cells.Item(138, 14).Value = -35359.712

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 12487.268
$ws.Cells.Item(32, 9).Value = 13983.432
$ws.Cells.Item(32, 10).Value = 7001.3335
$ws.Cells.Item(32, 11).Value = 13983.432
$ws.Cells.Item(32, 12).Value = 7001.3335
$ws.Cells.Item(32, 13).Value = -13696.432
$ws.Cells.Item(32, 14).Value = -7575.3335

$ws.Cells.Item(45, 8).Value = 1143.4286
$ws.Cells.Item(45, 9).Value = 1084
$ws.Cells.Item(45, 10).Value = 1500
$ws.Cells.Item(45, 11).Value = 1084
$ws.Cells.Item(45, 12).Value = 1500
$ws.Cells.Item(45, 13).Value = -707
$ws.Cells.Item(45, 14).Value = -2254

$ws.Cells.Item(110, 8).Value = 916.82855
$ws.Cells.Item(110, 9).Value = 940.7857
$ws.Cells.Item(110, 10).Value = 821
$ws.Cells.Item(110, 11).Value = 940.7857
$ws.Cells.Item(110, 12).Value = 821
$ws.Cells.Item(110, 13).Value = 1104.2143
$ws.Cells.Item(110, 14).Value = -4911

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 14).ClearContents()

$ws.Cells.Item(107, 8).Value = 65116.625
$ws.Cells.Item(107, 9).Value = 65116.625
$ws.Cells.Item(107, 10).Value = 0
$ws.Cells.Item(107, 11).Value = 65116.625
$ws.Cells.Item(107, 12).Value = 0
$ws.Cells.Item(107, 13).Value = -63196.625
$ws.Cells.Item(107, 14).ClearContents()

$ws.Cells.Item(134, 8).Value = 1781.6279
$ws.Cells.Item(134, 9).Value = 1614.575
$ws.Cells.Item(134, 10).Value = 4009
$ws.Cells.Item(134, 11).Value = 4843.725
$ws.Cells.Item(134, 12).Value = 12027
$ws.Cells.Item(134, 13).Value = -2308.725
$ws.Cells.Item(134, 14).Value = -17097

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2497.0435
$ws.Cells.Item(31, 9).Value = 1671.6
$ws.Cells.Item(31, 10).Value = 8000
$ws.Cells.Item(31, 11).Value = 1671.6
$ws.Cells.Item(31, 12).Value = 8000
$ws.Cells.Item(31, 13).Value = -1376.6
$ws.Cells.Item(31, 14).Value = -8590

$ws.Cells.Item(34, 8).Value = 2497.0435
$ws.Cells.Item(34, 9).Value = 1671.6
$ws.Cells.Item(34, 10).Value = 8000
$ws.Cells.Item(34, 11).Value = 1671.6
$ws.Cells.Item(34, 12).Value = 8000
$ws.Cells.Item(34, 13).Value = -1469.6
$ws.Cells.Item(34, 14).Value = -8404

$ws.Cells.Item(58, 8).Value = 1090891
$ws.Cells.Item(58, 9).Value = 1950856.5
$ws.Cells.Item(58, 10).Value = 1601.4
$ws.Cells.Item(58, 11).Value = 1950856.5
$ws.Cells.Item(58, 12).Value = 1601.4
$ws.Cells.Item(58, 13).Value = -1950653.5
$ws.Cells.Item(58, 14).Value = -2007.4

$ws.Cells.Item(99, 8).Value = 3640.1875
$ws.Cells.Item(99, 9).Value = 3928.818
$ws.Cells.Item(99, 10).Value = 3005.2
$ws.Cells.Item(99, 11).Value = 3928.818
$ws.Cells.Item(99, 12).Value = 3005.2
$ws.Cells.Item(99, 13).Value = -2430.818
$ws.Cells.Item(99, 14).Value = -6001.2

$ws.Cells.Item(126, 8).Value = 3640.1875
$ws.Cells.Item(126, 9).Value = 3928.818
$ws.Cells.Item(126, 10).Value = 3005.2
$ws.Cells.Item(126, 11).Value = 11786.454
$ws.Cells.Item(126, 12).Value = 9015.599999999999
$ws.Cells.Item(126, 13).Value = -9316.454000000002
$ws.Cells.Item(126, 14).Value = -13955.6

$ws.Cells.Item(127, 8).Value = 0
$ws.Cells.Item(127, 9).Value = 0
$ws.Cells.Item(127, 10).Value = 0
$ws.Cells.Item(127, 11).Value = 0
$ws.Cells.Item(127, 12).Value = 0
$ws.Cells.Item(127, 14).ClearContents()

$ws.Cells.Item(132, 8).Value = 266223.22
$ws.Cells.Item(132, 9).Value = 366245.53
$ws.Cells.Item(132, 10).Value = 1878.5
$ws.Cells.Item(132, 11).Value = 1098736.59
$ws.Cells.Item(132, 12).Value = 5635.5
$ws.Cells.Item(132, 13).Value = -1096206.59
$ws.Cells.Item(132, 14).Value = -10695.5

$ws.Cells.Item(134, 8).Value = 1135.6464
$ws.Cells.Item(134, 9).Value = 1018.89703
$ws.Cells.Item(134, 10).Value = 1702.7142
$ws.Cells.Item(134, 11).Value = 3056.69109
$ws.Cells.Item(134, 12).Value = 5108.142599999999
$ws.Cells.Item(134, 13).Value = -521.6910899999998
$ws.Cells.Item(134, 14).Value = -10178.1426

$ws.Cells.Item(136, 8).Value = 1090891
$ws.Cells.Item(136, 9).Value = 1950856.5
$ws.Cells.Item(136, 10).Value = 1601.4
$ws.Cells.Item(136, 11).Value = 5852569.5
$ws.Cells.Item(136, 12).Value = 4804.200000000001
$ws.Cells.Item(136, 13).Value = -5850019.5
$ws.Cells.Item(136, 14).Value = -9904.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(64, 8).Value = 4811.0557
$ws.Cells.Item(64, 9).Value = 1600
$ws.Cells.Item(64, 10).Value = 4999.9414
$ws.Cells.Item(64, 11).Value = 4800
$ws.Cells.Item(64, 12).Value = 14999.8242
$ws.Cells.Item(64, 13).Value = -4530
$ws.Cells.Item(64, 14).Value = -15539.8242

$ws.Cells.Item(67, 8).Value = 4811.0557
$ws.Cells.Item(67, 9).Value = 1600
$ws.Cells.Item(67, 10).Value = 4999.9414
$ws.Cells.Item(67, 11).Value = 4800
$ws.Cells.Item(67, 12).Value = 14999.8242
$ws.Cells.Item(67, 13).Value = -3864
$ws.Cells.Item(67, 14).Value = -16871.8242

$ws.Cells.Item(106, 8).Value = 6323.3335
$ws.Cells.Item(106, 9).Value = 0
$ws.Cells.Item(106, 10).Value = 6323.3335
$ws.Cells.Item(106, 11).Value = 0
$ws.Cells.Item(106, 12).Value = 18970.0005
$ws.Cells.Item(106, 14).Value = -20862.0005

$ws.Cells.Item(129, 8).Value = 1667946.9
$ws.Cells.Item(129, 9).Value = 667.3333
$ws.Cells.Item(129, 10).Value = 2779466.5
$ws.Cells.Item(129, 11).Value = 2001.9999
$ws.Cells.Item(129, 12).Value = 8338399.5
$ws.Cells.Item(129, 13).Value = 2998.0001
$ws.Cells.Item(129, 14).Value = -8348399.5

$ws.Cells.Item(131, 8).Value = 16952600
$ws.Cells.Item(131, 9).Value = 22186
$ws.Cells.Item(131, 10).Value = 18520230
$ws.Cells.Item(131, 11).Value = 66558
$ws.Cells.Item(131, 12).Value = 55560690
$ws.Cells.Item(131, 13).Value = -61518
$ws.Cells.Item(131, 14).Value = -55570770

$ws.Cells.Item(138, 8).Value = 1373.1818
$ws.Cells.Item(138, 9).Value = 817.375
$ws.Cells.Item(138, 10).Value = 2855.3333
$ws.Cells.Item(138, 11).Value = 2452.125
$ws.Cells.Item(138, 12).Value = 8565.999899999999
$ws.Cells.Item(138, 13).Value = 2687.875
$ws.Cells.Item(138, 14).Value = -18845.9999

$ws.Cells.Item(139, 8).Value = 2109.8147
$ws.Cells.Item(139, 9).Value = 1810.4736
$ws.Cells.Item(139, 10).Value = 2820.75
$ws.Cells.Item(139, 11).Value = 5431.4208
$ws.Cells.Item(139, 12).Value = 8462.25
$ws.Cells.Item(139, 13).Value = -291.4207999999999
$ws.Cells.Item(139, 14).Value = -18742.25

$ws.Cells.Item(140, 8).Value = 2202.9312
$ws.Cells.Item(140, 9).Value = 1448.5238
$ws.Cells.Item(140, 10).Value = 4183.25
$ws.Cells.Item(140, 11).Value = 4345.5714
$ws.Cells.Item(140, 12).Value = 12549.75
$ws.Cells.Item(140, 13).Value = 834.4286000000002
$ws.Cells.Item(140, 14).Value = -22909.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 1385.2903
$ws.Cells.Item(132, 9).Value = 941.9167
$ws.Cells.Item(132, 10).Value = 2905.4285
$ws.Cells.Item(132, 11).Value = 2825.7501
$ws.Cells.Item(132, 12).Value = 8716.2855
$ws.Cells.Item(132, 13).Value = -295.7501000000002
$ws.Cells.Item(132, 14).Value = -13776.2855

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(35, 8).Value = 2842.8572
$ws.Cells.Item(35, 9).Value = 2483.3333
$ws.Cells.Item(35, 10).Value = 5000
$ws.Cells.Item(35, 11).Value = 2483.3333
$ws.Cells.Item(35, 12).Value = 5000
$ws.Cells.Item(35, 13).Value = -2147.3333
$ws.Cells.Item(35, 14).Value = -5672

$ws.Cells.Item(40, 8).Value = 4098.0625
$ws.Cells.Item(40, 9).Value = 6226.5
$ws.Cells.Item(40, 10).Value = 2821
$ws.Cells.Item(40, 11).Value = 6226.5
$ws.Cells.Item(40, 12).Value = 2821
$ws.Cells.Item(40, 13).Value = -6090.5
$ws.Cells.Item(40, 14).Value = -3093

$ws.Cells.Item(93, 8).Value = 1042.9286
$ws.Cells.Item(93, 9).Value = 775.25
$ws.Cells.Item(93, 10).Value = 1399.8334
$ws.Cells.Item(93, 11).Value = 775.25
$ws.Cells.Item(93, 12).Value = 1399.8334
$ws.Cells.Item(93, 13).Value = 472.75
$ws.Cells.Item(93, 14).Value = -3895.8334

$ws.Cells.Item(117, 8).Value = 44999
$ws.Cells.Item(117, 9).Value = 0
$ws.Cells.Item(117, 10).Value = 44999
$ws.Cells.Item(117, 11).Value = 0
$ws.Cells.Item(117, 12).Value = 44999
$ws.Cells.Item(117, 14).Value = -54177

$ws.Cells.Item(132, 8).Value = 2595.6338
$ws.Cells.Item(132, 9).Value = 2136.4482
$ws.Cells.Item(132, 10).Value = 4644.3076
$ws.Cells.Item(132, 11).Value = 6409.344599999999
$ws.Cells.Item(132, 12).Value = 13932.9228
$ws.Cells.Item(132, 13).Value = -3879.344599999999
$ws.Cells.Item(132, 14).Value = -18992.9228

$ws.Cells.Item(136, 8).Value = 3526.2856
$ws.Cells.Item(136, 9).Value = 3848.575
$ws.Cells.Item(136, 10).Value = 2720.5625
$ws.Cells.Item(136, 11).Value = 11545.725
$ws.Cells.Item(136, 12).Value = 8161.6875
$ws.Cells.Item(136, 13).Value = -8995.724999999999
$ws.Cells.Item(136, 14).Value = -13261.6875

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(57, 8).Value = 45365
$ws.Cells.Item(57, 9).Value = 52000
$ws.Cells.Item(57, 10).Value = 42047.5
$ws.Cells.Item(57, 11).Value = 52000
$ws.Cells.Item(57, 12).Value = 42047.5
$ws.Cells.Item(57, 13).Value = -51246
$ws.Cells.Item(57, 14).Value = -43555.5

$ws.Cells.Item(107, 8).Value = 552.24243
$ws.Cells.Item(107, 9).Value = 506.7037
$ws.Cells.Item(107, 10).Value = 757.1667
$ws.Cells.Item(107, 11).Value = 1520.1111
$ws.Cells.Item(107, 12).Value = 2271.5001
$ws.Cells.Item(107, 13).Value = 399.8888999999999
$ws.Cells.Item(107, 14).Value = -6111.5001

$ws.Cells.Item(118, 8).Value = 24139.2
$ws.Cells.Item(118, 9).Value = 0
$ws.Cells.Item(118, 10).Value = 24139.2
$ws.Cells.Item(118, 11).Value = 0
$ws.Cells.Item(118, 12).Value = 24139.2
$ws.Cells.Item(118, 14).Value = -27453.2
